# Insert two new weekly-report rows (Fecha = 44722, 2022-06-10) right above the
# existing row for Fecha = 44281 (the previously-oldest "Primera"/"Segunda" pair),
# shifting all the rows below down by two. This mirrors the commit
# "Fruta / hortaliza, semanal": a new week's prices are prepended to the
# historical series kept in this sheet, pushing everything else down by 2 rows
# and extending the used range from A1:R593 to A1:R595.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 489 downwards (2 whole rows inserted before the old row 489).
$insertRange = $ws.Range("A489:A490").EntireRow
$insertRange.Insert()

# --- New row 489: Acelga, "Primera", Fecha 44722 ---
$ws.Cells.Item(489, 1).Value  = 9
$ws.Cells.Item(489, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(489, 3).Value  = "Metropolitana"
$ws.Cells.Item(489, 4).Value  = 44722
$ws.Cells.Item(489, 5).Value  = 13
$ws.Cells.Item(489, 6).Value  = 100112009
$ws.Cells.Item(489, 7).Value  = "Acelga"
$ws.Cells.Item(489, 8).Value  = "Sin especificar"
$ws.Cells.Item(489, 9).Value  = "Primera"
$ws.Cells.Item(489, 10).Value = 52
$ws.Cells.Item(489, 11).Value = 11000
$ws.Cells.Item(489, 12).Value = 11000
$ws.Cells.Item(489, 13).Value = 11000
$ws.Cells.Item(489, 14).Value = "$/docena de atados"
$ws.Cells.Item(489, 15).Value = "Región Metropolitana"
$ws.Cells.Item(489, 16).Value = 3667
$ws.Cells.Item(489, 17).Value = 3
$ws.Cells.Item(489, 18).Value = "Hortaliza"

# --- New row 490: Acelga, "Segunda", Fecha 44722 ---
$ws.Cells.Item(490, 1).Value  = 9
$ws.Cells.Item(490, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(490, 3).Value  = "Metropolitana"
$ws.Cells.Item(490, 4).Value  = 44722
$ws.Cells.Item(490, 5).Value  = 13
$ws.Cells.Item(490, 6).Value  = 100112009
$ws.Cells.Item(490, 7).Value  = "Acelga"
$ws.Cells.Item(490, 8).Value  = "Sin especificar"
$ws.Cells.Item(490, 9).Value  = "Segunda"
$ws.Cells.Item(490, 10).Value = 34
$ws.Cells.Item(490, 11).Value = 9000
$ws.Cells.Item(490, 12).Value = 9000
$ws.Cells.Item(490, 13).Value = 9000
$ws.Cells.Item(490, 14).Value = "$/docena de atados"
$ws.Cells.Item(490, 15).Value = "Región Metropolitana"
$ws.Cells.Item(490, 16).Value = 3000
$ws.Cells.Item(490, 17).Value = 3
$ws.Cells.Item(490, 18).Value = "Hortaliza"
